# This script updates the "Price" (D) and "Volume(1h)" (E) columns of the
# cryptocurrency tracking sheet to reflect refreshed values from a scheduled
# GitHub Actions run, per the upstream commit
# "Updated cryptos list on Sat Feb 18 20:53:08 UTC 2023 with GitHub Actions".
#
# All of these cells store plain text (prices/percentages are kept as
# strings, e.g. "1.002" or "  -0.03%  ", not numeric values) in the source
# workbook. Excel's COM Value setter auto-detects and coerces plain
# numeric-looking strings (like "1.002" or "52.90") into real numbers,
# which would both change the cell type and silently drop formatting such
# as trailing zeros. To prevent that, Set-TextCell below temporarily marks
# the target cell as Text ("@") before assigning the value, then calls
# ClearFormats() to drop the now-unneeded number-format override so the
# cell's style index is left exactly as it was originally (unstyled).
# Values that can never be parsed as a number (e.g. "24.675.14", which has
# two '.' separators, or the percentage strings that include spaces/%)
# are assigned directly since Excel already keeps them as text.
#
# NOTE: parameters are passed positionally (not with -Name) for
# compatibility with this PowerShell host's function-binding behavior.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

function Set-PlainCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $ws.Range($Address).Value = $Text
}

Set-PlainCell "D2" '24.675.14'
Set-PlainCell "E2" '  -0.03%  '
Set-PlainCell "D3" '1.687.31'
Set-PlainCell "E3" '  -0.88%  '
Set-TextCell "D4" '1.002'
Set-PlainCell "E4" '  +0.58%  '
Set-PlainCell "E5" '  +0.15%  '
Set-PlainCell "E6" '  +0.65%  '
Set-TextCell "D7" '0.3934'
Set-PlainCell "E7" '  -1.08%  '
Set-TextCell "D8" '0.4037'
Set-PlainCell "E8" '  -0.72%  '
Set-PlainCell "E9" '  +0.55%  '
Set-TextCell "D10" '1.482'
Set-PlainCell "E10" '  -2.57%  '
Set-TextCell "D11" '52.90'
Set-PlainCell "E11" '  -0.76%  '
Set-TextCell "D12" '0.08802'
Set-PlainCell "E12" '  +0.18%  '
Set-TextCell "D13" '7.233'
Set-PlainCell "E13" '  -1.43%  '
Set-PlainCell "E14" '  +0.43%  '
Set-TextCell "D15" '8.041'
Set-PlainCell "E16" '  -1.05%  '
Set-PlainCell "D17" '1.693.99'
Set-PlainCell "E17" '  -0.43%  '
Set-PlainCell "E18" '  -1.67%  '
Set-TextCell "D19" '0.07007'
Set-PlainCell "E19" '  -1.40%  '
Set-PlainCell "E20" '  -0.34%  '
Set-TextCell "D21" '6.978'
Set-PlainCell "E21" '  +3.37%  '
Set-PlainCell "E22" '  +0.96%  '
Set-PlainCell "E23" '  +0.38%  '
Set-PlainCell "D24" '24.650.00'
Set-PlainCell "E24" '  -0.08%  '
Set-TextCell "D25" '3.288'
Set-PlainCell "E25" '  +9.94%  '
Set-TextCell "D26" '2.362'
Set-PlainCell "E26" '  +2.36%  '
Set-TextCell "D27" '22.67'
Set-PlainCell "E27" '  +1.10%  '
Set-TextCell "D28" '162.45'
Set-PlainCell "E28" '  +2.12%  '
Set-TextCell "D29" '5.189'
Set-PlainCell "E29" '  +1.10%  '
Set-TextCell "D30" '135.06'
Set-PlainCell "E30" '  +1.28%  '
Set-TextCell "D31" '7.561'
Set-PlainCell "E31" '  +1.14%  '
Set-PlainCell "D32" '1.878.95'
Set-PlainCell "E32" '  -0.39%  '
Set-PlainCell "E33" '  -3.33%  '
Set-TextCell "D34" '0.08527'
Set-PlainCell "E34" '  -1.78%  '
Set-TextCell "D35" '7.128'
Set-PlainCell "E35" '  -3.61%  '
Set-PlainCell "E36" '  +0.04%  '
Set-PlainCell "E37" '  -0.31%  '
Set-TextCell "D38" '1.879'
Set-PlainCell "E38" '  -3.43%  '
Set-PlainCell "E39" '  -3.24%  '
Set-TextCell "D40" '0.09156'
Set-PlainCell "E40" '  +1.75%  '
Set-TextCell "D41" '0.02704'
Set-PlainCell "E41" '  -2.34%  '
Set-TextCell "D42" '1.460'
Set-PlainCell "E42" '  -1.04%  '
Set-TextCell "D43" '0.7573'
Set-PlainCell "E43" '  -1.09%  '
Set-TextCell "D44" '15.94'
Set-PlainCell "E44" '  +2.57%  '
Set-TextCell "D45" '2.586'
Set-PlainCell "E45" '  +5.19%  '
Set-TextCell "D46" '0.7105'
Set-PlainCell "E46" '  -1.86%  '
Set-TextCell "D47" '4.221'
Set-PlainCell "E47" '  +1.24%  '
Set-PlainCell "E48" '  +0.64%  '
Set-TextCell "D49" '139.35'
Set-PlainCell "E49" '  -1.53%  '
Set-TextCell "D50" '1.307'
Set-PlainCell "E50" '  -1.10%  '
Set-TextCell "D51" '0.07956'
Set-PlainCell "E51" '  -0.77%  '

Write-Host "Applied cryptos list update (85 cell changes across D2:E51)."
